$wb = $excel.ActiveWorkbook

# Rename Sheet1 to AddCustomerTest
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "AddCustomerTest"

# Populate header row
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "postCode"

# Populate data row
$ws.Range("A2").Value = "sudhir"
$ws.Range("B2").Value = "chakravarthi"
$ws.Range("C2").Value = 515001

# Match selection / active cell from the diff
$ws.Range("C2").Select()
